$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row for new "Launch Tests" data block
$ws.Range("I1").Value = "Launch_Angle"
$ws.Range("J1").Value = "Acceleration"
$ws.Range("K1").Value = "Range"

# Row 2
$ws.Range("I2").Value = 45
$ws.Range("J2").Value = 3000
$ws.Range("K2").Value = 2.5
$ws.Range("L2").Value = 2.5
$ws.Range("M2").Value = 2.5
$ws.Range("N2").Value = 2.5
$ws.Range("O2").Value = 2.5

# Row 3
$ws.Range("I3").Value = 35
$ws.Range("J3").Value = 5000
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 2
$ws.Range("M3").Value = 2
$ws.Range("N3").Value = 2
$ws.Range("O3").Value = 2

# Row 4
$ws.Range("I4").Value = 45
$ws.Range("J4").Value = 5000
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 2
$ws.Range("M4").Value = 2
$ws.Range("N4").Value = 2
$ws.Range("O4").Value = 2

# Row 5
$ws.Range("I5").Value = 45
$ws.Range("J5").Value = 2750
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 2
$ws.Range("M5").Value = 2
$ws.Range("N5").Value = 2
$ws.Range("O5").Value = 2

# Row 6
$ws.Range("I6").Value = 60
$ws.Range("J6").Value = 2750
$ws.Range("K6").Value = 1.5
$ws.Range("L6").Value = 1.5
$ws.Range("M6").Value = 1.5
$ws.Range("N6").Value = 1.5
$ws.Range("O6").Value = 1.5

# Row 7
$ws.Range("I7").Value = 50
$ws.Range("J7").Value = 2750
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 2
$ws.Range("M7").Value = 2
$ws.Range("N7").Value = 2
$ws.Range("O7").Value = 2

# Update selection to match target (E22) as noted in the diff
$ws.Range("E22").Select()
